# Macroferia Regional de Talca - Tomate: weekly refresh.
# Two new daily observations were inserted ahead of the existing series
# (pushing rows 594-627 down to 596-629), and the two columns of figures
# for the new week are populated into the freshly inserted rows 594-595.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 594, shifting the rest
# of the table (594..627) down to (596..629). Excel carries the existing
# row formatting (e.g. the date style on column D) down onto the new rows.
$ws.Rows.Item(594).Resize(2).Insert()

# New row 594
$ws.Range("A594").Value = 5
$ws.Range("B594").Value = "Macroferia Regional de Talca"
$ws.Range("C594").Value = "Maule"
$ws.Range("D594").Value = 44706
$ws.Range("E594").Value = 7
$ws.Range("F594").Value = 100112020
$ws.Range("G594").Value = "Tomate"
$ws.Range("H594").Value = "Larga vida"
$ws.Range("I594").Value = "Primera"
$ws.Range("J594").Value = 2500
$ws.Range("K594").Value = 17000
$ws.Range("L594").Value = 17000
$ws.Range("M594").Value = 17000
$ws.Range("N594").Value = "$/bandeja 18 kilos"
$ws.Range("O594").Value = "Región de Arica y Parinacota"
$ws.Range("P594").Value = 944
$ws.Range("Q594").Value = 18
$ws.Range("R594").Value = "Hortaliza"

# New row 595
$ws.Range("A595").Value = 5
$ws.Range("B595").Value = "Macroferia Regional de Talca"
$ws.Range("C595").Value = "Maule"
$ws.Range("D595").Value = 44706
$ws.Range("E595").Value = 7
$ws.Range("F595").Value = 100112020
$ws.Range("G595").Value = "Tomate"
$ws.Range("H595").Value = "Larga vida"
$ws.Range("I595").Value = "Primera"
$ws.Range("J595").Value = 2500
$ws.Range("K595").Value = 8000
$ws.Range("L595").Value = 8000
$ws.Range("M595").Value = 8000
$ws.Range("N595").Value = "$/caja 10 kilos"
$ws.Range("O595").Value = "Región de Arica y Parinacota"
$ws.Range("P595").Value = 800
$ws.Range("Q595").Value = 10
$ws.Range("R595").Value = "Hortaliza"
